$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("Sheet1"))
$newSheet.Name = "DataSet2"
$newSheet.Range("B2").Value = "b@b.b"
$newSheet.Hyperlinks.Add($newSheet.Range("B2"), "mailto:b@b.b")
